# Applies the commit: split the single "ODI Batting" sheet into three sheets:
#   1) "Player Info"       - new small player metadata table
#   2) "ODI Batting"       - the original batting log, MATCH_CARD_LINK -> MATCH_CODE
#                            (URL collapsed to just the numeric match code) plus
#                            two newly scraped rows
#   3) "ODI Batting Extra" - new supplementary per-match batting stats table

$wb = $excel.ActiveWorkbook
$orig = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Step 1: build the new "ODI Batting" sheet (sheetId will become 2) by copying
# the current sheet's data (header + 139 rows) as values, then re-apply the
# header formatting (bold/border/center) that lives on the original header.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$batting = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)

$orig.Range("A1:J139").Copy()
$batting.Range("A1:J139").PasteSpecial(-4163)   # xlPasteValues
$orig.Range("A1:J1").Copy()
$batting.Range("A1:J1").PasteSpecial(-4122)     # xlPasteFormats

# Rename the MATCH_CARD_LINK header to MATCH_CODE
$batting.Range("D1").Value = "MATCH_CODE"

# Collapse every MATCH_CARD_LINK url in column D down to the bare match code
for ($r = 2; $r -le 139; $r++) {
    $cell = $batting.Cells.Item($r, 4)
    $url = $cell.Value2
    $parts = $url -split "MatchCode="
    $code = $parts[1]
    $cell.Value = "'" + $code
}

# Append the two newly scraped rows (140, 141), forcing every value to be
# stored as text (matches the source data's inlineStr typing)
$newRows = @(
    @("139","139","31/03/2023","4746","2nd","Netherlands","Willowmoore Park","c F J Klaassen b Aryan Dutt","9","21"),
    @("140","140","02/04/2023","4751","1st","Netherlands","Wanderers Stadium","c T L W Cooper b V J Kingma","8","19")
)
$rowNum = 140
foreach ($rowData in $newRows) {
    for ($c = 1; $c -le 10; $c++) {
        $val = $rowData[$c - 1]
        $batting.Cells.Item($rowNum, $c).Value = "'" + $val
    }
    $rowNum++
}

# ---------------------------------------------------------------------------
# Step 2: build the new "ODI Batting Extra" sheet (sheetId will become 3)
# ---------------------------------------------------------------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet2)

$extraHeaders = @("MATCH_CODE","BATTING_POSITION","NUM_4","NUM_6","PERCENT_RUNS_OF_TOTAL","MAN_OF_MATCH")
for ($c = 1; $c -le 6; $c++) {
    $extra.Cells.Item(1, $c).Value = $extraHeaders[$c - 1]
}
$orig.Range("A1:F1").Copy()
$extra.Range("A1:F1").PasteSpecial(-4122)       # xlPasteFormats

# rows: MATCH_CODE, BATTING_POSITION(number|blank), NUM_4(text|blank), NUM_6(text|blank), PERCENT_RUNS_OF_TOTAL(text|blank), MAN_OF_MATCH(text)
$extraData = @(
    @("4421", $null, $null, $null, $null, "NO"),
    @("4458", 2, "3", "1", "7.33%", "NO"),
    @("4459", 1, "10", "1", "23.46%", "NO"),
    @("4478", $null, $null, $null, $null, "NO"),
    @("4524", $null, $null, $null, $null, "NO"),
    @("4526", 2, "7", "3", "27.08%", "YES"),
    @("4529", 1, "12", "2", "43.21%", "YES"),
    @("4557", 2, "9", "2", "31.79%", "NO"),
    @("4559", 2, "2", "0", "7.79%", "NO"),
    @("4619", $null, $null, $null, $null, "NO"),
    @("4620", 1, "1", "0", "6.02%", "NO"),
    @("4622", $null, $null, $null, $null, "NO"),
    @("4656", $null, $null, $null, $null, "NO"),
    @("4657", 1, "1", "0", "1.80%", "NO"),
    @("4658", 2, "1", "0", "6.06%", "NO"),
    @("4698", 1, "5", "1", "12.42%", "NO"),
    @("4699", 1, "4", "1", "8.93%", "NO"),
    @("4727", $null, $null, $null, $null, "NO"),
    @("4746", 1, "2", "0", "4.74%", "NO"),
    @("4751", $null, $null, $null, $null, "NO")
)
$rowNum = 2
foreach ($rowData in $extraData) {
    $extra.Cells.Item($rowNum, 1).Value = "'" + $rowData[0]
    if ($null -ne $rowData[1]) {
        $extra.Cells.Item($rowNum, 2).Value = $rowData[1]
    }
    if ($null -ne $rowData[2]) {
        $extra.Cells.Item($rowNum, 3).Value = "'" + $rowData[2]
    }
    if ($null -ne $rowData[3]) {
        $extra.Cells.Item($rowNum, 4).Value = "'" + $rowData[3]
    }
    if ($null -ne $rowData[4]) {
        $extra.Cells.Item($rowNum, 5).Value = "'" + $rowData[4]
    }
    $extra.Cells.Item($rowNum, 6).Value = $rowData[5]
    $rowNum++
}

# ---------------------------------------------------------------------------
# Step 3: turn the original sheet into "Player Info" (keeps sheetId 1)
# ---------------------------------------------------------------------------
$orig.Cells.Clear()

$playerHeaders = @("ID","NAME","BATTING_HAND","BOWL_STYLE")
for ($c = 1; $c -le 4; $c++) {
    $orig.Cells.Item(1, $c).Value = $playerHeaders[$c - 1]
}
$batting.Range("A1:D1").Copy()
$orig.Range("A1:D1").PasteSpecial(-4122)        # xlPasteFormats (reuse header style)

$orig.Cells.Item(2, 1).Value = "'3988"
$orig.Cells.Item(2, 2).Value = "Quinton de Kock"
$orig.Cells.Item(2, 3).Value = "Left Handed"
$orig.Cells.Item(2, 4).Value = "Does Not Bowl | Unknown"

# ---------------------------------------------------------------------------
# Step 4: rename sheets now that names are free of collisions, in final order
# ---------------------------------------------------------------------------
$orig.Name = "Player Info"
$batting.Name = "ODI Batting"
$extra.Name = "ODI Batting Extra"
